$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update today's conversion rates in the note ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 3.41 = 13106.35 pesos
✅ 13106.35 pesos = 3.41 = 955.67 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%"

# --- tasas sheet: updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 293.499
$ws2.Range("O10").Value = 3846.7
$ws2.Range("N12").Value = 3840
$ws2.Range("O12").Value = 280
